$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.06490000000001
$ws.Range("C4").Value = -11.2033
$ws.Range("C5").Value = -14.60070000000001
$ws.Range("A7").Value = -21.59610000000001
$ws.Range("C8").Value = -11.73250000000001
$ws.Range("A16").Value = -20.16939999999999
$ws.Range("C16").Value = -11.84030000000001
